$d = $word.ActiveDocument

# Change 1: "Below is a description/notes..." paragraph -> split into 4 runs:
# "Below are" + " " + "descriptions" + "/notes on the approaches I took in the execution of this assignment."
$para1 = $d.Paragraphs.Item(11)
$rng1 = $para1.Range
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00FA0381" w:rsidRDefault="009A3AB6"><w:pPr><w:rPr><w:lang w:val="en-IE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IE"/></w:rPr><w:t>Below are</w:t></w:r><w:r><w:rPr><w:lang w:val="en-IE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-IE"/></w:rPr><w:t>descriptions</w:t></w:r><w:r><w:rPr><w:lang w:val="en-IE"/></w:rPr><w:t>/notes on the approaches I took in the execution of this assignment.</w:t></w:r></w:p>'
$rng1.InsertXML($xml1)

# Change 2: "Used this resource for formats" + break + hyperlink -> single sentence about
# division-by-zero check, removing the hyperlink, keeping the trailing line break.
$para2 = $d.Paragraphs.Item(14)
$rng2 = $para2.Range
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009A3AB6" w:rsidRDefault="00277B82" w:rsidP="00FA0381"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-IE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-IE"/></w:rPr><w:t>For all arithmetic division calculation a test was done to check that the divisor was non-zero.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-IE"/></w:rPr><w:br/></w:r></w:p>'
$rng2.InsertXML($xml2)
